$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: new supplier record
$ws.Range("A2").Value = 3001
$ws.Range("B2").Value = "MADI INTERNATIONNAL"
$ws.Range("C2").Value = "MADI INTERNATIONNAL"
$ws.Range("D2").Value = 5005001

# Column N (14) width, as left behind by the paste/formatting that introduced the new row
$ws.Columns.Item(14).ColumnWidth = 20.6

# Final selection left on screen after the edit
$ws.Range("A3:D10").Select()
